$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1).Range
$p1.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ReleCloud</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Announces the Launch of a New Delivery Drone</w:t></w:r></w:p>
'@)

$p2 = $d.Paragraphs(2).Range
$p2.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Subtitle"/></w:pPr><w:r><w:t>The innovative drone will offer faster, safer, and greener delivery options for customers</w:t></w:r></w:p>
'@)

$p3 = $d.Paragraphs(3).Range
$p3.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ReleCloud</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, a leading cloud computing and e-commerce company, has announced the launch of a new delivery drone that will revolutionize the online shopping experience. The drone, named </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ReleFly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, will be available for customers in summer 2024 and will offer faster, safer, and greener delivery options.</w:t></w:r></w:p>
'@)

$p4 = $d.Paragraphs(4).Range
$p4.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ReleFly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is a lightweight, autonomous, and electric drone that can carry packages up to 5 kg and fly up to 25 km in a single charge.</w:t></w:r></w:p>
'@)

$p5 = $d.Paragraphs(5).Range
$p5.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ReleFly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> uses advanced artificial intelligence and computer vision to navigate complex urban environments and avoid obstacles. It can also communicate with other drones and traffic systems to ensure safe and efficient delivery.</w:t></w:r></w:p>
'@)

$p6 = $d.Paragraphs(6).Range
$p6.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ReleFly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is designed to reduce carbon emissions and noise pollution by using renewable energy and silent propellers. It can also land and take off vertically, minimizing the need for infrastructure and space.</w:t></w:r></w:p>
'@)

$p7 = $d.Paragraphs(7).Range
$p7.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:lastRenderedPageBreak/><w:t>ReleCloud's</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> CEO, Jane Smith, said that </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ReleFly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is a game-changer for the e-commerce industry and the customers. "We are excited to introduce </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ReleFly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, our new delivery drone, that will offer our customers a faster, safer, and greener way to receive their orders. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ReleFly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is the result of years of research and development and reflects our commitment to innovation and customer satisfaction. With </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ReleFly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, we aim to make online shopping more convenient and enjoyable for everyone."</w:t></w:r></w:p>
'@)

$p8 = $d.Paragraphs(8).Range
$p8.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ReleCloud</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is a global cloud computing and e-commerce company that offers a wide range of products and services, including online retail, cloud computing, digital streaming, artificial intelligence, and more.</w:t></w:r></w:p>
'@)

$p9 = $d.Paragraphs(9).Range
$p9.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ReleCloud</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> has over 200 million customers worldwide and employs over 500,000 people. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ReleCloud's</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> revenue in 2023 was $500 billion, making it one of the most valuable companies in the world.</w:t></w:r></w:p>
'@)

$p10 = $d.Paragraphs(10).Range
$p10.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ReleCloud</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is known for its customer-centric culture and its focus on innovation and social responsibility. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ReleCloud</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> has invested heavily in renewable energy, education, health, and philanthropy.</w:t></w:r></w:p>
'@)

$p11 = $d.Paragraphs(11).Range
$p11.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t xml:space="preserve">For more information about </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ReleFly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ReleCloud</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, please visit www.relecloud.com or contact press@relecloud.com.</w:t></w:r></w:p>
'@)

